$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.873.58'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.049.17'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''245.14'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = '''0.654'
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').Value = '''57.50'
$ws.Range('E7').Value = '  -2.93%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '''58.92'
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('D11').Value = '''0.0777'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('D13').Value = '''15.10'
$ws.Range('E13').Value = '  -4.40%  '
$ws.Range('D14').Value = '''0.873'
$ws.Range('E14').Value = '  +4.96%  '
$ws.Range('D15').Value = '2.347.38'
$ws.Range('E15').Value = '  +0.12%  '
$ws.Range('D16').Value = '''5.56'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '2.003.65'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '36.824.89'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '''17.41'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '''73.09'
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '0.0₃0888'
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('D22').Value = '''5.41'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').Value = '''236.01'
$ws.Range('E23').Value = '  -0.67%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = '''10.00'
$ws.Range('E26').Value = '  +6.66%  '
$ws.Range('D27').Value = '''2.20'
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').Value = '''168.73'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '''20.10'
$ws.Range('E29').Value = '  +0.27%  '
$ws.Range('D30').Value = '''5.46'
$ws.Range('E30').Value = '  +14.39%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').Value = '''1.14'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('E33').Value = '  +5.87%  '
$ws.Range('D34').Value = '''0.0615'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('D35').Value = '''2.37'
$ws.Range('E35').Value = '  +7.14%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +5.60%  '
$ws.Range('D38').Value = '''0.0845'
$ws.Range('E38').Value = '  -5.99%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('D40').Value = '''0.0223'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  -6.68%  '
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('D44').Value = '''0.0955'
$ws.Range('E44').Value = '  -10.33%  '
$ws.Range('D45').Value = '''96.67'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').Value = '''16.75'
$ws.Range('E46').Value = '  -4.11%  '
$ws.Range('D47').Value = '1.306.75'
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('D48').Value = '''2.36'
$ws.Range('E48').Value = '  -4.32%  '
$ws.Range('D49').Value = '''2.85'
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').Value = '''6.76'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = '2.234.28'
$ws.Range('E51').Value = '  +0.17%  '

Write-Output "Applied 89 cell updates"
